$d = $word.ActiveDocument

# 1. MIU 1308 -> MIU 1309
$d.Content.Find.Execute("MIU 1308", $false, $false, $false, $false, $false, $true, 1, $false, "MIU 1309", 2)

# 2. Remove the hyperlink (replace with a plain, unstyled run)
$h = $d.Hyperlinks.Item(1)
$hStart = $h.Range.Start
$displayText = $h.TextToDisplay
$h.Delete()
$hRange2 = $d.Range($hStart, $hStart + $displayText.Length)
$hRange2.Delete()
$insertPoint = $d.Range($hStart, $hStart)
$insertPoint.InsertAfter($displayText)

# 3. Move the _GoBack bookmark from paragraph 3 to the final (merged) empty paragraph
$d.Bookmarks.Item("_GoBack").Delete()

# 4. Collapse the two trailing empty paragraphs into a single one by deleting
#    the paragraph mark that separates them
$count = $d.Paragraphs.Count
$secondLastPara = $d.Paragraphs.Item($count - 1)
$markRange = $d.Range($secondLastPara.Range.End - 1, $secondLastPara.Range.End)
$markRange.Delete()

# 5. Re-insert the _GoBack bookmark into the (now single) trailing empty paragraph
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $finalPara.Range)
